# Fix the timetable header cell: "10-4" was being misread (Excel's auto-
# correction / reopening turned the hyphen into something unexpected), so
# re-enter the class/group label using a dot instead of a dash: "10.4".
#
# A1 (merged A1:A2) already carries a Text number format (numFmtId 49) with
# center/center alignment, so simply writing the corrected string keeps it
# stored as text (no accidental reinterpretation as a date/number) and
# preserves all existing formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "10.4"
